$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for credit card types (L1:O1)
$ws.Range("L1").Value = "mastercard"
$ws.Range("M1").Value = "visa"
$ws.Range("N1").Value = "discovercard"
$ws.Range("O1").Value = "americanexpress"

# Update the selected cell on the sheet to match the target state
$ws.Range("L8").Select()
